$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.286.12'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.673.27'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.72'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '194.95'
$ws.Range('E6').Value = '  +7.80%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.74'
$ws.Range('E10').Value = '  +4.68%  '
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('E12').Value = '  -4.87%  '
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.260.30'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.673.45'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.95'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.108.16'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '404.03'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.43'
$ws.Range('E22').Value = '  -3.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '88.18'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.09'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.69'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.07'
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('E28').Value = '  -3.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.40'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.99'
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.28'
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.37'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '44.97'
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '66.16'
$ws.Range('E34').Value = '  +3.88%  '
$ws.Range('E35').Value = '  +0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '612.23'
$ws.Range('E36').Value = '  +4.04%  '
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.397'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  -11.00%  '
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0429'
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.54'
$ws.Range('E44').Value = '  -6.39%  '
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.805.58'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.21'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.99'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '143.83'
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.59'
$ws.Range('E50').Value = '  -3.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.56'
$ws.Range('E51').Value = '  -10.84%  '
